# Agrego algo de contenido al archivo Word
$d = $word.ActiveDocument

# The document currently has a single paragraph with the text "Holaaa".
# We need to turn it into two runs: "Holaaa " (note trailing space) and
# "Como estás??".
$para = $d.Paragraphs.First
$r = $para.Range
$r.Collapse(0)  # wdCollapseEnd -> collapse to the end of the paragraph's range (before the pilcrow)
$r.InsertAfter(" Como estás??")
